$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header rotation: J<-K, K<-L, L<-M, M<-J(original)
$ws.Range("J1").Value = 'Sub-ontology'
$ws.Range("K1").Value = 'Definition source'
$ws.Range("L1").Value = 'Cross reference'
$ws.Range("M1").Value = 'Informal label for repository'

# Row 6: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J6").Value = ''
$ws.Range("K6").Value = 'Based on https://www.ncbi.nlm.nih.gov/pmc/articles/PMC6858509/'
$ws.Range("L6").Value = ''
$ws.Range("M6").Value = ""

# Row 7: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J7").Value = 'Intervention content and delivery'
$ws.Range("K7").Value = ''
$ws.Range("L7").Value = ''
$ws.Range("M7").Value = ""

# Row 12: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J12").Value = ''
$ws.Range("K12").Value = 'Insipred by https://www.merriam-webster.com/dictionary/bodybuilding'
$ws.Range("L12").Value = ''
$ws.Range("M12").Value = ""

# Row 14: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J14").Value = ''
$ws.Range("K14").Value = 'Cross-reference: COPPER:1044 - changed parent class'
$ws.Range("L14").Value = 'COPPER:1044'
$ws.Range("M14").Value = ""

# Row 19: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J19").Value = 'Intervention content and delivery'
$ws.Range("K19").Value = ''
$ws.Range("L19").Value = ''
$ws.Range("M19").Value = ""

# Row 20: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J20").Value = ''
$ws.Range("K20").Value = 'Parent class from Ontology for Biomedical Investigations
'
$ws.Range("L20").Value = ''
$ws.Range("M20").Value = ""

# Row 21: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J21").Value = 'Intervention content and delivery'
$ws.Range("K21").Value = ''
$ws.Range("L21").Value = ''
$ws.Range("M21").Value = ""

# Row 31: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J31").Value = 'intervention content and delivery'
$ws.Range("K31").Value = 'inspired by BCIO:008525'
$ws.Range("L31").Value = 'BCIO:008525'
$ws.Range("M31").Value = ""

# Row 32: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J32").Value = ''
$ws.Range("K32").Value = 'Based on http://humanbehaviourchange.org/ontology/BCIO_007000
'
$ws.Range("L32").Value = 'http://humanbehaviourchange.org/ontology/BCIO_007000'
$ws.Range("M32").Value = ""

# Row 34: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J34").Value = 'intervention content and delivery'
$ws.Range("K34").Value = 'inspired by BCIO:008560
'
$ws.Range("L34").Value = 'BCIO:008560'
$ws.Range("M34").Value = ""

# Row 36: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J36").Value = ''
$ws.Range("K36").Value = 'Drew on: https://www.oxfordlearnersdictionaries.com/definition/english/jogging#:~:text=jogging-,noun,as%20a%20form%20of%20exercise'
$ws.Range("L36").Value = ''
$ws.Range("M36").Value = ""

# Row 46: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J46").Value = 'Intervention content and delivery'
$ws.Range("K46").Value = ''
$ws.Range("L46").Value = ''
$ws.Range("M46").Value = ""

# Row 49: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J49").Value = 'Intervention content and delivery'
$ws.Range("K49").Value = ''
$ws.Range("L49").Value = ''
$ws.Range("M49").Value = ""

# Row 50: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J50").Value = 'Intervention content and delivery'
$ws.Range("K50").Value = 'BCIO:036042; GMHO:0000239'
$ws.Range("L50").Value = ''
$ws.Range("M50").Value = ""

# Row 59: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J59").Value = ''
$ws.Range("K59").Value = 'https://bciosearch.org/BCIO_050364'
$ws.Range("L59").Value = 'BCIO:050364'
$ws.Range("M59").Value = ""

# Row 61: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J61").Value = ''
$ws.Range("K61").Value = 'Crosss-reference: COPPER:1005 - changed parent class - definition in their ontology adapted definitions from NCIT_C154219 '
$ws.Range("L61").Value = ' COPPER:1005'
$ws.Range("M61").Value = ""

# Row 62: shift J<-K, K<-L, L<-M, M<-empty
$ws.Range("J62").Value = ''
$ws.Range("K62").Value = '"Rubenson, J., Heliams, D. B., Lloyd, D. G., and Fournier, P. A. (2004). Gait selection in the ostrich: mechanical and metabolic characteristics of walking and running with and without an aerial phase. Proceedings of the Royal Society of London. Series B: Biological Sciences, 271(1543), 1091-1099.
Also drew on:
- COPPER:1011"'
$ws.Range("L62").Value = ''
$ws.Range("M62").Value = ""

# Row 8: K -> J
$ws.Range("J8").Value = 'Intervention content and delivery'
$ws.Range("K8").Value = ""

# Row 17: K -> J
$ws.Range("J17").Value = 'Intervention content and delivery'
$ws.Range("K17").Value = ""

# Row 28: K -> J
$ws.Range("J28").Value = 'Intervention content and delivery'
$ws.Range("K28").Value = ""

# Row 41: J -> M
$ws.Range("M41").Value = 'unit'
$ws.Range("J41").Value = ""

# Row 47: normalize style to Normal (remove fill), restructure columns, and update Curation status
$ws.Range("A47:W47").Style = "Normal"
$ws.Range("J47").Value = 'Intervention content and delivery'
$ws.Range("K47").Value = ""
$ws.Range("T47").Value = 'External'